$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 453.18
$ws.Range("I15").Value = 453.18
$ws.Range("K15").Value = 1359.54
$ws.Range("M15").Value = -1190.54
$ws.Range("H116").Value = 327147.06
$ws.Range("I116").Value = 4694.4116
$ws.Range("K116").Value = 4694.4116
$ws.Range("M116").Value = -1252.4116
$ws.Range("H137").Value = 1666.4445
$ws.Range("I137").Value = 1584.1428
$ws.Range("J137").Value = 1718.8182
$ws.Range("K137").Value = 4752.428400000001
$ws.Range("L137").Value = 5156.4546
$ws.Range("M137").Value = -2202.428400000001
$ws.Range("N137").Value = -10256.4546
$ws.Range("H141").Value = 7851.6665
$ws.Range("I141").Value = 7090
$ws.Range("K141").Value = 21270
$ws.Range("M141").Value = -16090

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9160.781000000001
$ws.Range("I32").Value = 5557.05
$ws.Range("K32").Value = 5557.05
$ws.Range("M32").Value = -5270.05
$ws.Range("H61").Value = 10562.4
$ws.Range("I61").Value = 16924.8
$ws.Range("J61").Value = 4200
$ws.Range("K61").Value = 16924.8
$ws.Range("L61").Value = 4200
$ws.Range("M61").Value = -16712.8
$ws.Range("N61").Value = -4624
$ws.Range("H74").Value = 1395.4375
$ws.Range("I74").Value = 929.9
$ws.Range("J74").Value = 2171.3333
$ws.Range("K74").Value = 929.9
$ws.Range("L74").Value = 2171.3333
$ws.Range("M74").Value = -55.89999999999998
$ws.Range("N74").Value = -3919.3333
$ws.Range("H77").Value = 1395.4375
$ws.Range("I77").Value = 929.9
$ws.Range("J77").Value = 2171.3333
$ws.Range("K77").Value = 4649.5
$ws.Range("L77").Value = 10856.6665
$ws.Range("M77").Value = -281.5
$ws.Range("N77").Value = -19592.6665
$ws.Range("H122").Value = 43479664
$ws.Range("I122").Value = 50001476
$ws.Range("J122").Value = 896.6667
$ws.Range("K122").Value = 150004428
$ws.Range("L122").Value = 2690.0001
$ws.Range("M122").Value = -150001978
$ws.Range("N122").Value = -7590.0001
$ws.Range("H132").Value = 2909126.2
$ws.Range("I132").Value = 5953760.5
$ws.Range("K132").Value = 17861281.5
$ws.Range("M132").Value = -17858751.5
$ws.Range("H136").Value = 10562.4
$ws.Range("I136").Value = 16924.8
$ws.Range("J136").Value = 4200
$ws.Range("K136").Value = 50774.39999999999
$ws.Range("L136").Value = 12600
$ws.Range("M136").Value = -48224.39999999999
$ws.Range("N136").Value = -17700

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 55640644
$ws.Range("I134").Value = 111278616
$ws.Range("J134").Value = 2671.3333
$ws.Range("K134").Value = 333835848
$ws.Range("L134").Value = 8013.999899999999
$ws.Range("M134").Value = -333833313
$ws.Range("N134").Value = -13083.9999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11206.516
$ws.Range("I31").Value = 6413.65
$ws.Range("J31").Value = 18580.154
$ws.Range("K31").Value = 6413.65
$ws.Range("L31").Value = 18580.154
$ws.Range("M31").Value = -6118.65
$ws.Range("N31").Value = -19170.154
$ws.Range("H34").Value = 11206.516
$ws.Range("I34").Value = 6413.65
$ws.Range("J34").Value = 18580.154
$ws.Range("K34").Value = 6413.65
$ws.Range("L34").Value = 18580.154
$ws.Range("M34").Value = -6211.65
$ws.Range("N34").Value = -18984.154
$ws.Range("H86").Value = 2984.125
$ws.Range("I86").Value = 2510
$ws.Range("J86").Value = 3774.3333
$ws.Range("K86").Value = 2510
$ws.Range("L86").Value = 3774.3333
$ws.Range("M86").Value = -1387
$ws.Range("N86").Value = -6020.3333
$ws.Range("H89").Value = 2984.125
$ws.Range("I89").Value = 2510
$ws.Range("J89").Value = 3774.3333
$ws.Range("K89").Value = 12550
$ws.Range("L89").Value = 18871.6665
$ws.Range("M89").Value = -6934
$ws.Range("N89").Value = -30103.6665
$ws.Range("H132").Value = 9808650
$ws.Range("I132").Value = 17545032
$ws.Range("J132").Value = 9233.467000000001
$ws.Range("K132").Value = 52635096
$ws.Range("L132").Value = 27700.401
$ws.Range("M132").Value = -52632566
$ws.Range("N132").Value = -32760.401
$ws.Range("H134").Value = 9767626
$ws.Range("I134").Value = 12501871
$ws.Range("K134").Value = 37505613
$ws.Range("M134").Value = -37503078

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 412.76923
$ws.Range("I5").Value = 270.41177
$ws.Range("J5").Value = 681.6667
$ws.Range("K5").Value = 811.23531
$ws.Range("L5").Value = 2045.0001
$ws.Range("M5").Value = -699.23531
$ws.Range("N5").Value = -2269.0001
$ws.Range("H122").Value = 576
$ws.Range("I122").Value = 459.0625
$ws.Range("J122").Value = 1199.6666
$ws.Range("K122").Value = 4131.5625
$ws.Range("L122").Value = 10796.9994
$ws.Range("M122").Value = -1681.5625
$ws.Range("N122").Value = -15696.9994
$ws.Range("H135").Value = 412.76923
$ws.Range("I135").Value = 270.41177
$ws.Range("J135").Value = 681.6667
$ws.Range("K135").Value = 2433.70593
$ws.Range("L135").Value = 6135.0003
$ws.Range("M135").Value = 101.2940699999999
$ws.Range("N135").Value = -11205.0003

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 33371718
$ws.Range("I132").Value = 71502216
$ws.Range("J132").Value = 7530.5625
$ws.Range("K132").Value = 214506648
$ws.Range("L132").Value = 22591.6875
$ws.Range("M132").Value = -214504118
$ws.Range("N132").Value = -27651.6875

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1505.0952
$ws.Range("I7").Value = 1520.35
$ws.Range("K7").Value = 1520.35
$ws.Range("M7").Value = -1408.35
$ws.Range("H126").Value = 1505.0952
$ws.Range("I126").Value = 1520.35
$ws.Range("K126").Value = 4561.049999999999
$ws.Range("M126").Value = -2091.049999999999
$ws.Range("H136").Value = 5242.9644
$ws.Range("I136").Value = 6068.316
$ws.Range("J136").Value = 3500.5557
$ws.Range("K136").Value = 18204.948
$ws.Range("L136").Value = 10501.6671
$ws.Range("M136").Value = -15654.948
$ws.Range("N136").Value = -15601.6671

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 50025000
$ws.Range("J63").Value = 50025000
$ws.Range("L63").Value = 50025000
$ws.Range("N63").Value = -50026248
$ws.Range("H66").Value = 50025000
$ws.Range("J66").Value = 50025000
$ws.Range("L66").Value = 150075000
$ws.Range("N66").Value = -150081240
$ws.Range("H132").Value = 13241166
$ws.Range("I132").Value = 5912524
$ws.Range("J132").Value = 25699858
$ws.Range("K132").Value = 17737572
$ws.Range("L132").Value = 77099574
$ws.Range("M132").Value = -17735042
$ws.Range("N132").Value = -77104634
